# Fruta / hortaliza, semanal
# Insert two new daily observation rows into the weekly logic sheet
# (Feria Lagunitas de Puerto Montt - Limón), pushing the existing
# rows 643.. down by two, and append the matching two rows that were
# previously at the tail of the range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 643-644; everything currently at row 643
# and below shifts down to 645 and below (dimension grows from
# A1:T739 to A1:T741).
$ws.Rows("643:644").Insert()

# New row 643
$ws.Cells.Item(643, 1).Value = 4
$ws.Cells.Item(643, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(643, 3).Value = "Los Lagos"
$ws.Cells.Item(643, 4).Value = 44984
$ws.Cells.Item(643, 5).Value = 10
$ws.Cells.Item(643, 6).Value = "Fruta"
$ws.Cells.Item(643, 7).Value = 100102
$ws.Cells.Item(643, 8).Value = "Cítricos"
$ws.Cells.Item(643, 9).Value = 100102003
$ws.Cells.Item(643, 10).Value = "Limón"
$ws.Cells.Item(643, 11).Value = "Sin especificar"
$ws.Cells.Item(643, 12).Value = "1a plateado"
$ws.Cells.Item(643, 13).Value = 600
$ws.Cells.Item(643, 14).Value = 28000
$ws.Cells.Item(643, 15).Value = 29000
$ws.Cells.Item(643, 16).Value = 28500
$ws.Cells.Item(643, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(643, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(643, 19).Value = 1583
$ws.Cells.Item(643, 20).Value = 18

# New row 644
$ws.Cells.Item(644, 1).Value = 4
$ws.Cells.Item(644, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(644, 3).Value = "Los Lagos"
$ws.Cells.Item(644, 4).Value = 44984
$ws.Cells.Item(644, 5).Value = 10
$ws.Cells.Item(644, 6).Value = "Fruta"
$ws.Cells.Item(644, 7).Value = 100102
$ws.Cells.Item(644, 8).Value = "Cítricos"
$ws.Cells.Item(644, 9).Value = 100102003
$ws.Cells.Item(644, 10).Value = "Limón"
$ws.Cells.Item(644, 11).Value = "Sin especificar"
$ws.Cells.Item(644, 12).Value = "2a plateado"
$ws.Cells.Item(644, 13).Value = 300
$ws.Cells.Item(644, 14).Value = 26000
$ws.Cells.Item(644, 15).Value = 26000
$ws.Cells.Item(644, 16).Value = 26000
$ws.Cells.Item(644, 17).Value = "`$/malla 18 kilos"
$ws.Cells.Item(644, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(644, 19).Value = 1444
$ws.Cells.Item(644, 20).Value = 18

# Ensure date formatting on the new D cells matches the rest of the
# Fecha column (same numeric date format as the surrounding rows).
$ws.Range("D643:D644").NumberFormat = $ws.Range("D645").NumberFormat
